# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (per commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 121
    $ws.Range("F5").Value = 1069
    $ws.Range("F7").Value = 2128
    $ws.Range("F10").Value = 4723
    $ws.Range("F16").Value = 161
    $ws.Range("F20").Value = 3643
    $ws.Range("F21").Value = 364
    $ws.Range("F22").Value = 591
    $ws.Range("F26").Value = 108
    $ws.Range("F27").Value = 16
    $ws.Range("F28").Value = 9
    $ws.Range("F30").Value = 220

    if ($sheetName -eq "展览") {
        $ws.Range("F33").Value = 808
        $ws.Range("F34").Value = 2263
        $ws.Range("F35").Value = 415
    } else {
        $ws.Range("F34").Value = 808
        $ws.Range("F35").Value = 2263
        $ws.Range("F36").Value = 415
    }
}

$wb.Save()
